# Updates the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.356.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.814.39'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '342.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3809'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3494'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.78'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('E10').Value = '  +4.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07730'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.005'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.96%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.626'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.238'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.812.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.57%  '
$ws.Range('E17').Value = '  +3.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06727'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '85.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.92%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E21').Value = '  +7.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.557'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '27.382.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.475'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.667'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.68%  '
$ws.Range('E27').Value = '  +14.58%  '
$ws.Range('E28').Value = '  +11.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '154.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.009.21'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '135.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.15%  '
$ws.Range('E32').Value = '  +6.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.028'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08773'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.712'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.608'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6977'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2273'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02417'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06464'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.941'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.300'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6516'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.72%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.042'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.178'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07337'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.45%  '
